$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187; this shifts the existing rows 187-248
# down to 188-249, preserving all of their data and formatting.
$ws.Rows(187).Insert()

# Populate the newly inserted row 187 with the new record
# (same categorical values as the old row 187, but a new date and new
# volume/price figures).
$ws.Cells.Item(187, 1).Value = 5
$ws.Cells.Item(187, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(187, 3).Value = "Maule"
$ws.Cells.Item(187, 4).Value = 44559
$ws.Cells.Item(187, 5).Value = 7
$ws.Cells.Item(187, 6).Value = 100112023
$ws.Cells.Item(187, 7).Value = "Brócoli"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 5000
$ws.Cells.Item(187, 11).Value = 500
$ws.Cells.Item(187, 12).Value = 500
$ws.Cells.Item(187, 13).Value = 500
$ws.Cells.Item(187, 14).Value = "$/unidad"
$ws.Cells.Item(187, 15).Value = "Región del Maule"
$ws.Cells.Item(187, 16).Value = 500
$ws.Cells.Item(187, 17).Value = 1
$ws.Cells.Item(187, 18).Value = "Hortaliza"
